$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosRegistro")
$ws.Range("C2").Value = "12@gmail.com"
